$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 (bold, bordered, centered) onto new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(10, 10),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
